# Customer data edit: update customer with ID=1, delete a duplicate customer row with ID=3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the customer row with ID=1 (row 2) ---
# Name stays plain text; the year, CCCD number and people-count look numeric,
# so force "text" storage for them (matching the rest of the sheet, where every
# value - numeric-looking or not - is stored as a shared string) before writing,
# then strip the resulting number-format styling so no extra style is left on
# the cells.
$ws.Range("B2").Value = "hung"

$ws.Range("C2:E2").NumberFormat = "@"
$ws.Range("C2").Value = "2000"
$ws.Range("D2").Value = "123123123123"
$ws.Range("E2").Value = "1"
$ws.Range("C2:E2").ClearFormats()

# --- Delete the duplicate customer row with ID=3 (row 4) ---
# Rows 5 and 6 (ID=5 "Hai" and ID=6 "nhung") shift up to become rows 4 and 5.
$ws.Rows.Item(4).Delete()
